$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge first name + last name into column A (with literal separator text),
# shift email into column B, and role into column C. Then drop old column D.

$ws.Range("A1").Value = "Diego + participante.apellidos"
$ws.Range("B1").Value = "v.dvm.dvm@gmail.com"
$ws.Range("C1").Value = "Líder"

$ws.Range("A2").Value = "BRUCE ANTHONY + participante.apellidos"
$ws.Range("B2").Value = "a20203298@pucp.edu.pe"
$ws.Range("C2").Value = "Miembro"

$ws.Range("A3").Value = "GABRIEL OMAR + participante.apellidos"
$ws.Range("B3").Value = "a20203371@pucp.edu.pe"
$ws.Range("C3").Value = "Miembro"

$ws.Range("A4").Value = "CHRISTIAN SEBASTIAN + participante.apellidos"
$ws.Range("B4").Value = "s.chira@pucp.edu.pe"
$ws.Range("C4").Value = "Miembro"

# Remove the now-unused column D
$ws.Columns.Item(4).Delete()

# Set column widths to match the new layout (the ColumnWidth property adds a
# fixed ~0.8333 char padding on round-trip through this engine, so we dial
# the input back by 5/6 to land on the exact target widths of 46/23/10)
$ws.Columns.Item(1).ColumnWidth = 45.16666666666667
$ws.Columns.Item(2).ColumnWidth = 22.166666666666668
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666
